# flash分配.xlsx - add "temper_value" (warm/cold select value) row to the
# global-variable table on Sheet1, just above the existing "resetbtcnt" row.
#
# Before: row 12 = resetbtcnt / u8 / 蓝牙重新连接次数 / 0x2F80
#         row 13 = zigbee_join_cnt / u8 / 烧写后第一次上电为配网 / 0x2F81
# After:  row 12 = temper_value / u8 / 冷暖选择值 / 0x2F0B   (new)
#         row 13 = resetbtcnt / u8 / 蓝牙重新连接次数 / 0x2F80
#         row 14 = zigbee_join_cnt / u8 / 烧写后第一次上电为配网 / 0x2F81

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 12-13 down to 13-14, leaving a blank row 12 to fill in.
$ws.Rows("12:12").Insert()

# Pick up the bordered table style from the row above so the new row
# matches the rest of the table (Insert leaves the new row unformatted).
$ws.Range("A11:E11").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new "temper_value" entry (serial number 11).
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "temper_value"
$ws.Range("C12").Value = "u8"
$ws.Range("D12").Value = "冷暖选择值"
$ws.Range("E12").Value = "0x2F0B"

# Renumber the serial numbers of the rows that shifted down.
$ws.Range("A13").Value = 12
$ws.Range("A14").Value = 13

# Match the author's final selection.
[void]$ws.Range("E13").Select()
